$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Change 1: merge the three runs "—Authored " + "the" + " book "
# (all sharing identical run formatting) into a single run reading
# "—Authored the book ".  A same-text Find/Replace on the exact
# existing wording re-serialises the three runs into one run
# carrying the shared formatting.
# ---------------------------------------------------------------
$d.Content.Find.Execute("—Authored the book ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "—Authored the book ", 2) | Out-Null

# ---------------------------------------------------------------
# Change 2: in the "—Grammarly" / "— Writing" achievements list,
# the paragraph with the hidden _GoBack bookmark changes from
#   [bookmark]"— Writing"
# to
#   "—"[bookmark]"Writing"
# i.e. a new run containing just the em-dash is inserted before the
# bookmark (with the same run formatting), and the existing run's
# text loses its "— " prefix, becoming just "Writing".
#
# Direct Range.InsertBefore() at the bookmark's start position does
# create a new run ahead of the bookmark, but a run created that way
# carries no rPr. To reproduce the formatting exactly, temporarily
# remove the bookmark, perform the plain text edit (so the existing
# run keeps its full/untouched rPr), then re-insert the bookmark at
# the boundary between the new "—" run and the "Writing" run.
# ---------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bmStart = $bm.Start
$bm.Delete()

# Remove the leading em-dash + space ("— ") from the run, leaving "Writing".
$prefix = $d.Range($bmStart, $bmStart + 2)
$prefix.Text = ""

# Insert a new run containing just the em-dash ahead of that text; it
# inherits the run formatting already present on the "Writing" run.
$insPoint = $d.Range($bmStart, $bmStart)
$insPoint.InsertBefore("—")

# Put the _GoBack bookmark back exactly between the two runs.
$d.Bookmarks.Add("_GoBack", $d.Range($bmStart + 1, $bmStart + 1)) | Out-Null
